$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.063.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.102.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "

$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.008"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5170"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4417"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09453"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "52.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.174"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "25.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.111.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.704"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.073"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.56%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001164"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "

$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06702"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("E21").Value = "  -0.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.182"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.156.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.334"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.360.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.59%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.536"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.154"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.39%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.242"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.960"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.158"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.55%  "

$ws.Range("E37").Value = "  -3.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02572"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.38%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06767"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2274"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6929"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.303"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.275"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.55%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.636"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.97%  "

$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000358"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.220"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.88%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07188"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.89%  "
